$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Cell = "D2"; Value = "42.641.43"; Numeric = $false },
    @{ Cell = "E2"; Value = "  +1.08%  "; Numeric = $false },
    @{ Cell = "D3"; Value = "2.303.09"; Numeric = $false },
    @{ Cell = "E3"; Value = "  +0.43%  "; Numeric = $false },
    @{ Cell = "E4"; Value = "  +0.17%  "; Numeric = $false },
    @{ Cell = "D5"; Value = "316.49"; Numeric = $true },
    @{ Cell = "E5"; Value = "  -0.56%  "; Numeric = $false },
    @{ Cell = "D6"; Value = "104.04"; Numeric = $true },
    @{ Cell = "E6"; Value = "  +0.56%  "; Numeric = $false },
    @{ Cell = "D7"; Value = "0.629"; Numeric = $true },
    @{ Cell = "E7"; Value = "  -0.16%  "; Numeric = $false },
    @{ Cell = "E8"; Value = "  +0.22%  "; Numeric = $false },
    @{ Cell = "D9"; Value = "0.606"; Numeric = $true },
    @{ Cell = "E9"; Value = "  +0.15%  "; Numeric = $false },
    @{ Cell = "D10"; Value = "39.85"; Numeric = $true },
    @{ Cell = "E10"; Value = "  +1.25%  "; Numeric = $false },
    @{ Cell = "D11"; Value = "0.0907"; Numeric = $true },
    @{ Cell = "E11"; Value = "  +0.31%  "; Numeric = $false },
    @{ Cell = "D12"; Value = "8.51"; Numeric = $true },
    @{ Cell = "E12"; Value = "  +3.05%  "; Numeric = $false },
    @{ Cell = "E13"; Value = "  +1.25%  "; Numeric = $false },
    @{ Cell = "D14"; Value = "0.989"; Numeric = $true },
    @{ Cell = "D15"; Value = "15.41"; Numeric = $true },
    @{ Cell = "E15"; Value = "  +1.18%  "; Numeric = $false },
    @{ Cell = "D16"; Value = "2.654.33"; Numeric = $false },
    @{ Cell = "E16"; Value = "  +0.60%  "; Numeric = $false },
    @{ Cell = "D17"; Value = "2.306.23"; Numeric = $false },
    @{ Cell = "E17"; Value = "  +0.50%  "; Numeric = $false },
    @{ Cell = "D18"; Value = "42.573.52"; Numeric = $false },
    @{ Cell = "E18"; Value = "  +1.40%  "; Numeric = $false },
    @{ Cell = "D19"; Value = "7.60"; Numeric = $true },
    @{ Cell = "E19"; Value = "  +2.29%  "; Numeric = $false },
    @{ Cell = "E20"; Value = "  +0.58%  "; Numeric = $false },
    @{ Cell = "D21"; Value = "13.31"; Numeric = $true },
    @{ Cell = "E21"; Value = "  +32.29%  "; Numeric = $false },
    @{ Cell = "D22"; Value = "73.94"; Numeric = $true },
    @{ Cell = "E22"; Value = "  +0.81%  "; Numeric = $false },
    @{ Cell = "D23"; Value = "3.53"; Numeric = $true },
    @{ Cell = "E23"; Value = "  -3.37%  "; Numeric = $false },
    @{ Cell = "D24"; Value = "268.49"; Numeric = $true },
    @{ Cell = "E24"; Value = "  -3.93%  "; Numeric = $false },
    @{ Cell = "D25"; Value = "2.23"; Numeric = $true },
    @{ Cell = "E25"; Value = "  -1.14%  "; Numeric = $false },
    @{ Cell = "D27"; Value = "10.89"; Numeric = $true },
    @{ Cell = "E27"; Value = "  +0.92%  "; Numeric = $false },
    @{ Cell = "E28"; Value = "  +0.08%  "; Numeric = $false },
    @{ Cell = "D29"; Value = "22.61"; Numeric = $true },
    @{ Cell = "E29"; Value = "  -1.74%  "; Numeric = $false },
    @{ Cell = "D30"; Value = "38.06"; Numeric = $true },
    @{ Cell = "E30"; Value = "  +6.21%  "; Numeric = $false },
    @{ Cell = "D31"; Value = "6.51"; Numeric = $true },
    @{ Cell = "E31"; Value = "  +11.56%  "; Numeric = $false },
    @{ Cell = "D32"; Value = "166.14"; Numeric = $true },
    @{ Cell = "E32"; Value = "  +1.90%  "; Numeric = $false },
    @{ Cell = "D33"; Value = "0.0883"; Numeric = $true },
    @{ Cell = "E33"; Value = "  +1.40%  "; Numeric = $false },
    @{ Cell = "D34"; Value = "2.71"; Numeric = $true },
    @{ Cell = "E34"; Value = "  -4.95%  "; Numeric = $false },
    @{ Cell = "E35"; Value = "  -2.89%  "; Numeric = $false },
    @{ Cell = "D36"; Value = "0.113"; Numeric = $true },
    @{ Cell = "E36"; Value = "  -0.69%  "; Numeric = $false },
    @{ Cell = "D37"; Value = "4.57"; Numeric = $true },
    @{ Cell = "E37"; Value = "  +1.19%  "; Numeric = $false },
    @{ Cell = "D38"; Value = "0.0353"; Numeric = $true },
    @{ Cell = "E38"; Value = "  +1.30%  "; Numeric = $false },
    @{ Cell = "D39"; Value = "2.78"; Numeric = $true },
    @{ Cell = "E39"; Value = "  -4.60%  "; Numeric = $false },
    @{ Cell = "D40"; Value = "3.69"; Numeric = $true },
    @{ Cell = "E40"; Value = "  -0.50%  "; Numeric = $false },
    @{ Cell = "D41"; Value = "1.63"; Numeric = $true },
    @{ Cell = "E41"; Value = "  +12.35%  "; Numeric = $false },
    @{ Cell = "D42"; Value = "99.13"; Numeric = $true },
    @{ Cell = "E42"; Value = "  -0.82%  "; Numeric = $false },
    @{ Cell = "D43"; Value = "70.16"; Numeric = $true },
    @{ Cell = "E43"; Value = "  +1.13%  "; Numeric = $false },
    @{ Cell = "B44"; Value = "FirstDigitalUSD"; Numeric = $false },
    @{ Cell = "C44"; Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"; Numeric = $false },
    @{ Cell = "D44"; Value = "1.00"; Numeric = $true },
    @{ Cell = "E44"; Value = "  +0.22%  "; Numeric = $false },
    @{ Cell = "B45"; Value = "Algorand"; Numeric = $false },
    @{ Cell = "C45"; Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"; Numeric = $false },
    @{ Cell = "D45"; Value = "0.225"; Numeric = $true },
    @{ Cell = "E45"; Value = "  +0.37%  "; Numeric = $false },
    @{ Cell = "D46"; Value = "117.29"; Numeric = $true },
    @{ Cell = "E46"; Value = "  +3.74%  "; Numeric = $false },
    @{ Cell = "D47"; Value = "12.31"; Numeric = $true },
    @{ Cell = "E47"; Value = "  +3.35%  "; Numeric = $false },
    @{ Cell = "D48"; Value = "80.32"; Numeric = $true },
    @{ Cell = "E48"; Value = "  +4.51%  "; Numeric = $false },
    @{ Cell = "D49"; Value = "1.643.72"; Numeric = $false },
    @{ Cell = "E49"; Value = "  +4.48%  "; Numeric = $false },
    @{ Cell = "D50"; Value = "5.30"; Numeric = $true },
    @{ Cell = "E50"; Value = "  +0.21%  "; Numeric = $false },
    @{ Cell = "D51"; Value = "8.87"; Numeric = $true },
    @{ Cell = "E51"; Value = "  -0.74%  "; Numeric = $false }
)

foreach ($change in $changes) {
    $cell = $ws.Range($change.Cell)
    if ($change.Numeric) {
        # Force text storage for values that otherwise look like numbers
        # (e.g. "316.49", "0.629") so they stay text cells like the source data.
        $cell.NumberFormat = "@"
        $cell.Value = $change.Value
        $cell.ClearFormats()
    } else {
        $cell.Value = $change.Value
    }
}

Write-Output "Applied $($changes.Count) cell updates"
